$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Inandout" label from B2 but keep its (red font) style
$ws.Range("B2").Value = $null

# New column E of reference numbers alongside existing rows 2-11
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 6
$ws.Range("E4").Value = 7
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 4
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 2
$ws.Range("E10").Value = 8
$ws.Range("E11").Value = 9

# Two additional reference rows appended below the table
$ws.Range("B12").Value = "sdcasdc"
$ws.Range("B13").Value = "sdcasddd"

# Update the active selection to reflect where the user ended up (next free row)
$ws.Range("B14").Select()
